$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire ID column (column A), shifting NAME/AGE/PROFESSION/LOCATION left
$ws.Range("A:A").Delete()

# Update the selected cell to match the post-edit state
$ws.Range("I13").Select()
